# Update the cryptos price/volume table (columns D and E, rows 2-51)
# on Sheet1 to reflect the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '29.455.65'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '1.877.53'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = "'0.7142"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = "'242.03"
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = "'0.9996"
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = "'0.3123"
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("D9").Value = "'0.07736"
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("D10").Value = "'25.08"
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("D11").Value = "'0.08398"
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("D12").Value = '1.888.54'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = "'5.255"
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = "'0.7193"
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = "'91.74"
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '29.462.21'
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").Value = "'0.000008199"
$ws.Range("E17").Value = '  +4.76%  '
$ws.Range("D18").Value = "'5.992"
$ws.Range("E18").Value = '  +2.41%  '
$ws.Range("D19").Value = "'244.84"
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").Value = '2.130.24'
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("D21").Value = "'13.23"
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").Value = "'7.950"
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").Value = "'0.1635"
$ws.Range("E25").Value = '  +2.31%  '
$ws.Range("D26").Value = "'163.87"
$ws.Range("E26").Value = '  +0.80%  '
$ws.Range("D27").Value = "'9.044"
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").Value = "'18.64"
$ws.Range("E28").Value = '  +2.09%  '
$ws.Range("D29").Value = "'1.510"
$ws.Range("E29").Value = '  +0.92%  '
$ws.Range("D30").Value = "'4.429"
$ws.Range("E30").Value = '  +0.91%  '
$ws.Range("D31").Value = "'1.299"
$ws.Range("E31").Value = '  -3.96%  '
$ws.Range("D32").Value = "'4.332"
$ws.Range("E32").Value = '  +5.80%  '
$ws.Range("D33").Value = "'0.05241"
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("D34").Value = "'1.932"
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").Value = "'0.7707"
$ws.Range("E35").Value = '  +6.81%  '
$ws.Range("D36").Value = "'1.177"
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("D37").Value = "'2.677"
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").Value = "'0.01868"
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("D39").Value = "'2.723"
$ws.Range("E39").Value = '  +1.14%  '
$ws.Range("D40").Value = '1.176.30'
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").Value = "'6.427"
$ws.Range("E41").Value = '  +5.16%  '
$ws.Range("D42").Value = "'73.79"
$ws.Range("E42").Value = '  +1.51%  '
$ws.Range("D43").Value = "'0.8916"
$ws.Range("E43").Value = '  -1.72%  '
$ws.Range("D44").Value = "'104.29"
$ws.Range("E44").Value = '  +2.15%  '
$ws.Range("D45").Value = "'0.9992"
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").Value = '2.027.44'
$ws.Range("E46").Value = '  +0.93%  '
$ws.Range("D47").Value = "'1.808"
$ws.Range("E47").Value = '  +0.95%  '
$ws.Range("D48").Value = "'0.5197"
$ws.Range("E48").Value = '  -1.74%  '
$ws.Range("D49").Value = "'9.432"
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("D50").Value = "'0.4326"
$ws.Range("E50").Value = '  +0.94%  '
$ws.Range("D51").Value = "'7.067"
$ws.Range("E51").Value = '  +0.51%  '
